$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target roster data (rows 2-19), reflecting the updated player list,
# positions and teams after the upload.
$data = @(
    @{ Row = 2;  Name = "Derrick White";             Pos = "PF";       Team = "Boston Celtics" },
    @{ Row = 3;  Name = "Duncan Robinson";           Pos = "SG,SF";    Team = "Miami Heat" },
    @{ Row = 4;  Name = "Julius Randle";             Pos = "PF";       Team = "Minnesota Timberwolves" },
    @{ Row = 5;  Name = "Cameron Johnson";           Pos = "SF,PF";    Team = "Brooklyn Nets" },
    @{ Row = 6;  Name = "Anthony Davis";             Pos = "PF,C";     Team = "Los Angeles Lakers" },
    @{ Row = 7;  Name = "Bam Adebayo";               Pos = "C";        Team = "Miami Heat" },
    @{ Row = 8;  Name = "Isaiah Joe";                Pos = "PG,SG";    Team = "Oklahoma City Thunder" },
    @{ Row = 9;  Name = "Kentavious Caldwell-Pope";  Pos = "SG,SF";    Team = "Orlando Magic" },
    @{ Row = 10; Name = "Isaiah Hartenstein";        Pos = "C";        Team = "Oklahoma City Thunder" },
    @{ Row = 11; Name = "Quentin Grimes";            Pos = "SG,SF";    Team = "Dallas Mavericks" },
    @{ Row = 12; Name = "Damian Lillard";            Pos = "PG";       Team = "Milwaukee Bucks" },
    @{ Row = 13; Name = "Harrison Barnes";           Pos = "SF,PF";    Team = "San Antonio Spurs" },
    @{ Row = 14; Name = "Ty Jerome";                 Pos = "PG,SG";    Team = "Cleveland Cavaliers" },
    @{ Row = 15; Name = "Brandon Miller";            Pos = "SG,SF";    Team = "Charlotte Hornets" },
    @{ Row = 16; Name = "Nick Richards";             Pos = "C";        Team = "Charlotte Hornets" },
    @{ Row = 17; Name = "LaMelo Ball";               Pos = "PG,SG";    Team = "Charlotte Hornets" },
    @{ Row = 18; Name = "Brandon Ingram";            Pos = "SG,SF,PF"; Team = "New Orleans Pelicans" },
    @{ Row = 19; Name = "Cade Cunningham";           Pos = "PG,SG";    Team = "Detroit Pistons" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Name
    $ws.Cells.Item($r, 2).Value = $item.Pos
    $ws.Cells.Item($r, 3).Value = $item.Team
}

$wb.Save()
